$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM update: populate "Have" (column C) and "Bought" (column E) quantities ---
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 5
$ws.Range("C7").Value = 2
$ws.Range("E7").Value = 10
$ws.Range("C8").Value = 4
$ws.Range("C12").Value = 3
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 2
$ws.Range("C21").Value = 1
$ws.Range("C22").Value = 3
$ws.Range("C23").Value = 5
$ws.Range("E24").Value = 20
$ws.Range("C25").Value = 4
$ws.Range("C26").Value = 2
$ws.Range("C27").Value = 3
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("C31").Value = 1

# --- Column width adjustments for DistributorPartNum1 / DistributorPartLink1 ---
$ws.Columns(15).ColumnWidth = 12.666666666666666
$ws.Columns(16).ColumnWidth = 144.33333333333331

# --- Sheet view: freeze panes at column G (7 columns) and scroll / select ---
$ws.Range("H1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C31").Select()
